# Update "想去人数" (want-to-go count) figures pulled from the latest
# scrape (output generated at 456a3b4) for the gh-pages data sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6898
$ws1.Range("F4").Value = 39
$ws1.Range("F8").Value = 64
$ws1.Range("F14").Value = 139
$ws1.Range("F19").Value = 5028
$ws1.Range("F22").Value = 467

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6898
$ws4.Range("F4").Value = 39
$ws4.Range("F5").Value = 450
$ws4.Range("F8").Value = 64
$ws4.Range("F11").Value = 15
$ws4.Range("F13").Value = 0
$ws4.Range("F14").Value = 139
$ws4.Range("F20").Value = 5028
$ws4.Range("F24").Value = 467
